$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 5.404689
$ws.Range("H2").Value = 16.214067
$ws.Range("I2").Value = 0.08747859311663772
$ws.Range("J2").Value = 0.09021076583983562
$ws.Range("M2").Value = 2.914938
$ws.Range("N2").Value = 8.744814
$ws.Range("O2").Value = 0.07105539873786189
$ws.Range("P2").Value = 0.07541359350096061
$ws.Range("Q2").Value = 15.754333344282
$ws.Range("R2").Value = 141.789000098538
$ws.Range("S2").Value = 0.006215826314929874
$ws.Range("T2").Value = 0.006803118024455706

$ws.Range("G3").Value = 5.404689
$ws.Range("H3").Value = 16.214067
$ws.Range("I3").Value = 0.08747859311663772
$ws.Range("J3").Value = 0.09021076583983562
$ws.Range("M3").Value = 7.339638666666666
$ws.Range("N3").Value = 22.018916
$ws.Range("O3").Value = 0.1789132228719201
$ws.Range("P3").Value = 0.1898868953137022
$ws.Range("Q3").Value = 39.668464365708
$ws.Range("R3").Value = 357.0161792913719
$ws.Range("S3").Value = 0.01565107702679902
$ws.Range("T3").Value = 0.01712984224919777

$ws.Range("G4").Value = 5.404689
$ws.Range("H4").Value = 16.214067
$ws.Range("I4").Value = 0.08747859311663772
$ws.Range("J4").Value = 0.09021076583983562
$ws.Range("M4").Value = 15.66105466666667
$ws.Range("N4").Value = 46.983164
$ws.Range("O4").Value = 0.3817585430617917
$ws.Range("P4").Value = 0.4051737671361525
$ws.Range("Q4").Value = 84.643129885332
$ws.Range("R4").Value = 761.788168967988
$ws.Range("S4").Value = 0.03339570025730289
$ws.Range("T4").Value = 0.03655103583156354

$ws.Range("G5").Value = 5.404689
$ws.Range("H5").Value = 16.214067
$ws.Range("I5").Value = 0.08747859311663772
$ws.Range("J5").Value = 0.09021076583983562
$ws.Range("M5").Value = 7.112307
$ws.Range("N5").Value = 14.224614
$ws.Range("O5").Value = 0.1733717183113625
$ws.Range("P5").Value = 0.1226703344295343
$ws.Range("Q5").Value = 38.43980740752301
$ws.Range("R5").Value = 230.638844445138
$ws.Range("S5").Value = 0.01516631400409201
$ws.Range("T5").Value = 0.01106618481471704

$ws.Range("G6").Value = 5.404689
$ws.Range("H6").Value = 16.214067
$ws.Range("I6").Value = 0.08747859311663772
$ws.Range("J6").Value = 0.09021076583983562
$ws.Range("M6").Value = 7.995517333333333
$ws.Range("N6").Value = 23.986552
$ws.Range("O6").Value = 0.1949011170170639
$ws.Range("P6").Value = 0.2068554096196504
$ws.Range("Q6").Value = 43.213284580776
$ws.Range("R6").Value = 388.919561226984
$ws.Range("S6").Value = 0.01704967551351393
$ws.Range("T6").Value = 0.01866058491990156

$ws.Range("I7").Value = 0.2448858138641327
$ws.Range("J7").Value = 0.2525342032254661
$ws.Range("M7").Value = 2.914938
$ws.Range("N7").Value = 8.744814
$ws.Range("O7").Value = 0.07105539873786189
$ws.Range("P7").Value = 0.07541359350096061
$ws.Range("Q7").Value = 44.102363852118
$ws.Range("R7").Value = 396.921274669062
$ws.Range("S7").Value = 0.01740045914936178
$ws.Range("T7").Value = 0.01904451174713427

$ws.Range("I8").Value = 0.2448858138641327
$ws.Range("J8").Value = 0.2525342032254661
$ws.Range("M8").Value = 7.339638666666666
$ws.Range("N8").Value = 22.018916
$ws.Range("O8").Value = 0.1789132228719201
$ws.Range("P8").Value = 0.1898868953137022
$ws.Range("Q8").Value = 111.0471011803364
$ws.Range("R8").Value = 999.4239106230278
$ws.Range("S8").Value = 0.04381331019404511
$ws.Range("T8").Value = 0.04795293581100327

$ws.Range("I9").Value = 0.2448858138641327
$ws.Range("J9").Value = 0.2525342032254661
$ws.Range("M9").Value = 15.66105466666667
$ws.Range("N9").Value = 46.983164
$ws.Range("O9").Value = 0.3817585430617917
$ws.Range("P9").Value = 0.4051737671361525
$ws.Range("Q9").Value = 236.9482751321791
$ws.Range("R9").Value = 2132.534476189612
$ws.Range("S9").Value = 0.09348725151727243
$ws.Range("T9").Value = 0.1023202344515888

$ws.Range("I10").Value = 0.2448858138641327
$ws.Range("J10").Value = 0.2525342032254661
$ws.Range("M10").Value = 7.112307
$ws.Range("N10").Value = 14.224614
$ws.Range("O10").Value = 0.1733717183113625
$ws.Range("P10").Value = 0.1226703344295343
$ws.Range("Q10").Value = 107.607623607077
$ws.Range("R10").Value = 645.6457416424621
$ws.Range("S10").Value = 0.04245627433970117
$ws.Range("T10").Value = 0.0309784551645639

$ws.Range("I11").Value = 0.2448858138641327
$ws.Range("J11").Value = 0.2525342032254661
$ws.Range("M11").Value = 7.995517333333333
$ws.Range("N11").Value = 23.986552
$ws.Range("O11").Value = 0.1949011170170639
$ws.Range("P11").Value = 0.2068554096196504
$ws.Range("Q11").Value = 120.9703995833129
$ws.Range("R11").Value = 1088.733596249816
$ws.Range("S11").Value = 0.04772851866375227
$ws.Range("T11").Value = 0.05223806605117583

$ws.Range("G12").Value = 17.564497
$ws.Range("H12").Value = 52.69349099999999
$ws.Range("I12").Value = 0.2842934138044583
$ws.Range("J12").Value = 0.2931725999334087
$ws.Range("M12").Value = 2.914938
$ws.Range("N12").Value = 8.744814
$ws.Range("O12").Value = 0.07105539873786189
$ws.Range("P12").Value = 0.07541359350096061
$ws.Range("Q12").Value = 51.199419756186
$ws.Range("R12").Value = 460.794777805674
$ws.Range("S12").Value = 0.02020058187642375
$ws.Range("T12").Value = 0.02210919927699784

$ws.Range("G13").Value = 17.564497
$ws.Range("H13").Value = 52.69349099999999
$ws.Range("I13").Value = 0.2842934138044583
$ws.Range("J13").Value = 0.2931725999334087
$ws.Range("M13").Value = 7.339638666666666
$ws.Range("N13").Value = 22.018916
$ws.Range("O13").Value = 0.1789132228719201
$ws.Range("P13").Value = 0.1898868953137022
$ws.Range("Q13").Value = 128.9170613417506
$ws.Range("R13").Value = 1160.253552075756
$ws.Range("S13").Value = 0.05086385090501604
$ws.Range("T13").Value = 0.05566963479240107

$ws.Range("G14").Value = 17.564497
$ws.Range("H14").Value = 52.69349099999999
$ws.Range("I14").Value = 0.2842934138044583
$ws.Range("J14").Value = 0.2931725999334087
$ws.Range("M14").Value = 15.66105466666667
$ws.Range("N14").Value = 46.983164
$ws.Range("O14").Value = 0.3817585430617917
$ws.Range("P14").Value = 0.4051737671361525
$ws.Range("Q14").Value = 275.0785477095027
$ws.Range("R14").Value = 2475.706929385524
$ws.Range("S14").Value = 0.108531439456053
$ws.Range("T14").Value = 0.1187858467361193

$ws.Range("G15").Value = 17.564497
$ws.Range("H15").Value = 52.69349099999999
$ws.Range("I15").Value = 0.2842934138044583
$ws.Range("J15").Value = 0.2931725999334087
$ws.Range("M15").Value = 7.112307
$ws.Range("N15").Value = 14.224614
$ws.Range("O15").Value = 0.1733717183113625
$ws.Range("P15").Value = 0.1226703344295343
$ws.Range("Q15").Value = 124.924094964579
$ws.Range("R15").Value = 749.544569787474
$ws.Range("S15").Value = 0.04928843765588215
$ws.Range("T15").Value = 0.0359635808794073

$ws.Range("G16").Value = 17.564497
$ws.Range("H16").Value = 52.69349099999999
$ws.Range("I16").Value = 0.2842934138044583
$ws.Range("J16").Value = 0.2931725999334087
$ws.Range("M16").Value = 7.995517333333333
$ws.Range("N16").Value = 23.986552
$ws.Range("O16").Value = 0.1949011170170639
$ws.Range("P16").Value = 0.2068554096196504
$ws.Range("Q16").Value = 140.4372402147813
$ws.Range("R16").Value = 1263.935161933032
$ws.Range("S16").Value = 0.05540910391108329
$ws.Range("T16").Value = 0.06064433824848315

$ws.Range("G17").Value = 5.613580499999999
$ws.Range("H17").Value = 11.227161
$ws.Range("I17").Value = 0.09085964511315853
$ws.Range("J17").Value = 0.06246494429911598
$ws.Range("M17").Value = 2.914938
$ws.Range("N17").Value = 8.744814
$ws.Range("O17").Value = 0.07105539873786189
$ws.Range("P17").Value = 0.07541359350096061
$ws.Range("Q17").Value = 16.363239115509
$ws.Range("R17").Value = 98.17943469305399
$ws.Range("S17").Value = 0.006456068312696103
$ws.Range("T17").Value = 0.004710705917433679

$ws.Range("G18").Value = 5.613580499999999
$ws.Range("H18").Value = 11.227161
$ws.Range("I18").Value = 0.09085964511315853
$ws.Range("J18").Value = 0.06246494429911598
$ws.Range("M18").Value = 7.339638666666666
$ws.Range("N18").Value = 22.018916
$ws.Range("O18").Value = 0.1789132228719201
$ws.Range("P18").Value = 0.1898868953137022
$ws.Range("Q18").Value = 41.20165249624599
$ws.Range("R18").Value = 247.2099149774759
$ws.Range("S18").Value = 0.0162559919361941
$ws.Range("T18").Value = 0.01186127433890247

$ws.Range("G19").Value = 5.613580499999999
$ws.Range("H19").Value = 11.227161
$ws.Range("I19").Value = 0.09085964511315853
$ws.Range("J19").Value = 0.06246494429911598
$ws.Range("M19").Value = 15.66105466666667
$ws.Range("N19").Value = 46.983164
$ws.Range("O19").Value = 0.3817585430617917
$ws.Range("P19").Value = 0.4051737671361525
$ws.Range("Q19").Value = 87.914591086234
$ws.Range("R19").Value = 527.487546517404
$ws.Range("S19").Value = 0.03468644574151084
$ws.Range("T19").Value = 0.02530915679562275

$ws.Range("G20").Value = 5.613580499999999
$ws.Range("H20").Value = 11.227161
$ws.Range("I20").Value = 0.09085964511315853
$ws.Range("J20").Value = 0.06246494429911598
$ws.Range("M20").Value = 7.112307
$ws.Range("N20").Value = 14.224614
$ws.Range("O20").Value = 0.1733717183113625
$ws.Range("P20").Value = 0.1226703344295343
$ws.Range("Q20").Value = 39.9255078852135
$ws.Range("R20").Value = 159.702031540854
$ws.Range("S20").Value = 0.01575249279842889
$ws.Range("T20").Value = 0.007662595607294788

$ws.Range("G21").Value = 5.613580499999999
$ws.Range("H21").Value = 11.227161
$ws.Range("I21").Value = 0.09085964511315853
$ws.Range("J21").Value = 0.06246494429911598
$ws.Range("M21").Value = 7.995517333333333
$ws.Range("N21").Value = 23.986552
$ws.Range("O21").Value = 0.1949011170170639
$ws.Range("P21").Value = 0.2068554096196504
$ws.Range("Q21").Value = 44.88348018981199
$ws.Range("R21").Value = 269.3008811388719
$ws.Range("S21").Value = 0.01770864632432861
$ws.Range("T21").Value = 0.01292121163986228

$ws.Range("G22").Value = 18.07044533333334
$ws.Range("H22").Value = 54.211336
$ws.Range("I22").Value = 0.2924825341016128
$ws.Range("J22").Value = 0.3016174867021735
$ws.Range("M22").Value = 2.914938
$ws.Range("N22").Value = 8.744814
$ws.Range("O22").Value = 0.07105539873786189
$ws.Range("P22").Value = 0.07541359350096061
$ws.Range("Q22").Value = 52.674227779056
$ws.Range("R22").Value = 474.068050011504
$ws.Range("S22").Value = 0.02078246308445039
$ws.Range("T22").Value = 0.02274605853493911

$ws.Range("G23").Value = 18.07044533333334
$ws.Range("H23").Value = 54.211336
$ws.Range("I23").Value = 0.2924825341016128
$ws.Range("J23").Value = 0.3016174867021735
$ws.Range("M23").Value = 7.339638666666666
$ws.Range("N23").Value = 22.018916
$ws.Range("O23").Value = 0.1789132228719201
$ws.Range("P23").Value = 0.1898868953137022
$ws.Range("Q23").Value = 132.6305392924195
$ws.Range("R23").Value = 1193.674853631776
$ws.Range("S23").Value = 0.05232899280986581
$ws.Range("T23").Value = 0.05727320812219758

$ws.Range("G24").Value = 18.07044533333334
$ws.Range("H24").Value = 54.211336
$ws.Range("I24").Value = 0.2924825341016128
$ws.Range("J24").Value = 0.3016174867021735
$ws.Range("M24").Value = 15.66105466666667
$ws.Range("N24").Value = 46.983164
$ws.Range("O24").Value = 0.3817585430617917
$ws.Range("P24").Value = 0.4051737671361525
$ws.Range("Q24").Value = 283.0022322163449
$ws.Range("R24").Value = 2547.020089947104
$ws.Range("S24").Value = 0.1116577060896525
$ws.Range("T24").Value = 0.122207493321258

$ws.Range("G25").Value = 18.07044533333334
$ws.Range("H25").Value = 54.211336
$ws.Range("I25").Value = 0.2924825341016128
$ws.Range("J25").Value = 0.3016174867021735
$ws.Range("M25").Value = 7.112307
$ws.Range("N25").Value = 14.224614
$ws.Range("O25").Value = 0.1733717183113625
$ws.Range("P25").Value = 0.1226703344295343
$ws.Range("Q25").Value = 128.522554837384
$ws.Range("R25").Value = 771.135329024304
$ws.Range("S25").Value = 0.05070819951325829
$ws.Range("T25").Value = 0.03699951796355123

$ws.Range("G26").Value = 18.07044533333334
$ws.Range("H26").Value = 54.211336
$ws.Range("I26").Value = 0.2924825341016128
$ws.Range("J26").Value = 0.3016174867021735
$ws.Range("M26").Value = 7.995517333333333
$ws.Range("N26").Value = 23.986552
$ws.Range("O26").Value = 0.1949011170170639
$ws.Range("P26").Value = 0.2068554096196504
$ws.Range("Q26").Value = 144.4825588837191
$ws.Range("R26").Value = 1300.343029953472
$ws.Range("S26").Value = 0.05700517260438582
$ws.Range("T26").Value = 0.06239120876022757
